# Revert file upload functionality
# Append a new data row (row 59) at the bottom of each of the four
# worksheets (DE_LFT_#1, DE_LFT_#2, DE_PLT_#1, DE_PLT_#2), mirroring the
# structure of the existing rows (time / length / ID / actual-length /
# checksum + their decimal counterparts).

$wb = $excel.ActiveWorkbook

$rows = @(
    @{
        Sheet = "DE_LFT_#1"
        A = 45845.43280092593
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
        D = "0x01,0x5C"
        E = "0x14"
        F = 380
        G = [double]"7.598631275147109e+23"
        H = 348
        I = 14
    },
    @{
        Sheet = "DE_LFT_#2"
        A = 45845.43280092593
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
        D = "0x01,0x60"
        E = "0xe"
        F = 380
        G = [double]"5.68432987514711e+23"
        H = 352
        I = 14
    },
    @{
        Sheet = "DE_PLT_#1"
        A = 45845.43280092593
        B = "0x00,0x82"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x00,0x7E"
        E = "0x7"
        F = 130
        G = [double]"5.68631262647114e+23"
        H = 126
        I = 7
    },
    @{
        Sheet = "DE_PLT_#2"
        A = 45845.43280092593
        B = "0x00,0x82"
        C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
        D = "0x00,0x7D"
        E = "0x3"
        F = 130
        G = [double]"9.85046333984776e+23"
        H = 125
        I = 3
    }
)

foreach ($rowData in $rows) {
    $ws = $wb.Worksheets.Item($rowData.Sheet)

    # New row is one past the current last used row (row 58 -> row 59).
    $usedRange = $ws.UsedRange
    $lastRow = $usedRange.Row + $usedRange.Rows.Count - 1
    $newRow = $lastRow + 1

    # Column A keeps the same date/time number format as the rows above it.
    $ws.Cells.Item($newRow, 1).NumberFormat = $ws.Cells.Item($newRow - 1, 1).NumberFormat
    $ws.Cells.Item($newRow, 1).Value = $rowData.A

    $ws.Cells.Item($newRow, 2).Value = $rowData.B
    $ws.Cells.Item($newRow, 3).Value = $rowData.C
    $ws.Cells.Item($newRow, 4).Value = $rowData.D
    $ws.Cells.Item($newRow, 5).Value = $rowData.E
    $ws.Cells.Item($newRow, 6).Value = $rowData.F
    $ws.Cells.Item($newRow, 7).Value = $rowData.G
    $ws.Cells.Item($newRow, 8).Value = $rowData.H
    $ws.Cells.Item($newRow, 9).Value = $rowData.I
}
